{"js": "async (context) => {\n  const replacements = [\n    [\"2024-04-24 Wednesday\", \"2024-04-25 Thursday\"],\n    [\"42-29=13\", \"36+20=56\"],\n    [\"81-57=24\", \"86-10=76\"],\n    [\"25+43=68\", \"47-28=19\"],\n    [\"63+28=91\", \"16+33=49\"],\n    [\"44-40=4\", \"75-52=23\"],\n    [\"32+29=61\", \"20+7=27\"],\n    [\"27+25=52\", \"48+48=96\"],\n    [\"55-38=17\", \"47-42=5\"],\n    [\"6-3=3\", \"48-32=16\"],\n    [\"38+26=64\", \"14+72=86\"],\n    [\"44+35=79\", \"1+4=5\"],\n    [\"32+54=86\", \"59-4=55\"],\n    [\"7+52=59\", \"51-23=28\"],\n    [\"34+64=98\", \"98-73=25\"],\n    [\"78-43=35\", \"78-36=42\"],\n    [\"82-36=46\", \"8-3=5\"],\n    [\"55+4=59\", \"19+21=40\"],\n    [\"62-27=35\", \"20+9=29\"],\n    [\"74-3=71\", \"49+35=84\"],\n    [\"93-80=13\", \"74-9=65\"],\n    [\"37+49=86\", \"2+56=58\"],\n    [\"47-13=34\", \"41-27=14\"],\n    [\"92-1=91\", \"35+15=50\"],\n    [\"93-16=77\", \"95-8=87\"],\n    [\"0+79=79\", \"36+8=44\"],\n    [\"27+61=88\", \"76+4=80\"],\n    [\"24+5=29\", \"96-85=11\"],\n    [\"93-65=28\", \"50+45=95\"],\n    [\"58-52=6\", \"49+8=57\"],\n    [\"12+54=66\", \"37-22=15\"],\n    [\"67+15=82\", \"70-57=13\"],\n    [\"79-67=12\", \"7+70=77\"],\n    [\"88-53=35\", \"8+87=95\"],\n    [\"13-10=3\", \"23+37=60\"],\n    [\"30+29=59\", \"62+16=78\"],\n    [\"14+21=35\", \"70+11=81\"],\n    [\"2+3=5\", \"42-4=38\"],\n    [\"47+7=54\", \"64-13=51\"],\n    [\"52-6=46\", \"39-6=33\"],\n    [\"54-53=1\", \"3-2=1\"],\n    [\"59+15=74\", \"54-51=3\"],\n    [\"37-13=24\", \"48+42=90\"],\n    [\"51+11=62\", \"77+21=98\"],\n    [\"90-81=9\", \"93-19=74\"],\n    [\"62-33=29\", \"86-63=23\"],\n    [\"51+39=90\", \"67-9=58\"],\n    [\"55-47=8\", \"15+41=56\"],\n    [\"16+49=65\", \"75+5=80\"],\n    [\"63-48=15\", \"83+3=86\"],\n    [\"70-29=41\", \"40+21=61\"],\n    [\"48+26=74\", \"34-33=1\"],\n    [\"52+39=91\", \"46-43=3\"],\n    [\"31+40=71\", \"69-55=14\"],\n    [\"34+46=80\", \"0+63=63\"],\n    [\"11+1=12\", \"45+47=92\"],\n    [\"80-12=68\", \"59-42=17\"],\n    [\"39+20=59\", \"18+79=97\"],\n    [\"8+28=36\", \"21+1=22\"],\n    [\"63+1=64\", \"45-4=41\"],\n    [\"23+18=41\", \"68-23=45\"],\n    [\"3+88=91\", \"30+64=94\"],\n    [\"28+52=80\", \"58-8=50\"],\n    [\"25+62=87\", \"9+47=56\"],\n    [\"94-75=19\", \"65-41=24\"],\n    [\"55+33=88\", \"17+46=63\"],\n    [\"30-26=4\", \"40+19=59\"],\n    [\"43+18=61\", \"85+3=88\"],\n    [\"3+34=37\", \"44+1=45\"],\n    [\"73-46=27\", \"81-36=45\"],\n    [\"6+42=48\", \"32+60=92\"],\n    [\"44+23=67\", \"46+2=48\"],\n    [\"17+59=76\", \"0+6=6\"],\n    [\"27-19=8\", \"18+80=98\"],\n    [\"89-71=18\", \"78-58=20\"],\n    [\"68-28=40\", \"60+16=76\"],\n    [\"85+6=91\", \"21+53=74\"],\n    [\"69-39=30\", \"94-57=37\"],\n    [\"34-29=5\", \"12+41=53\"],\n    [\"3+92=95\", \"38+34=72\"],\n    [\"89-70=19\", \"98-0=98\"],\n    [\"82-4=78\", \"51+34=85\"],\n    [\"75-45=30\", \"94-93=1\"],\n    [\"66-49=17\", \"78-45=33\"],\n    [\"34+21=55\", \"70-18=52\"],\n    [\"30-29=1\", \"71-17=54\"],\n    [\"13-2=11\", \"9+62=71\"],\n    [\"38+50=88\", \"10-4=6\"],\n    [\"34+26=60\", \"44+32=76\"],\n    [\"54+40=94\", \"8+83=91\"],\n    [\"79+20=99\", \"33+55=88\"],\n    [\"85-32=53\", \"69-7=62\"],\n    [\"96-86=10\", \"13+85=98\"],\n    [\"91-70=21\", \"69-43=26\"],\n    [\"10+38=48\", \"12+14=26\"],\n    [\"71+19=90\", \"89-18=71\"],\n    [\"46-42=4\", \"5+51=56\"],\n    [\"81+8=89\", \"92-80=12\"],\n    [\"73-5=68\", \"90+6=96\"],\n    [\"88-30=58\", \"65+7=72\"],\n    [\"34-31=3\", \"92-27=65\"],\n  ];\n\n  const body = context.document.body;\n  const paras = body.paragraphs;\n  paras.load(\"items/text\");\n  await context.sync();\n\n  const items = paras.items;\n  if (items.length !== replacements.length) {\n    throw new Error(\"Paragraph count mismatch: \" + items.length + \" vs \" + replacements.length);\n  }\n\n  for (let i = 0; i < items.length; i++) {\n    const [oldText, newText] = replacements[i];\n    if (items[i].text !== oldText) {\n      throw new Error(\"Mismatch at paragraph \" + i + \": expected '\" + oldText + \"' but found '\" + items[i].text + \"'\");\n    }\n    items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-ExactText($doc, $oldText, $newText) {\n  $find = $doc.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = $wdFindContinue\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.MatchSoundsLike = $false\n  $find.MatchAllWordForms = $false\n  $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n}\n\nReplace-ExactText $d \"2024-04-24 Wednesday\" \"2024-04-25 Thursday\"\nReplace-ExactText $d \"42-29=13\" \"36+20=56\"\nReplace-ExactText $d \"81-57=24\" \"86-10=76\"\nReplace-ExactText $d \"25+43=68\" \"47-28=19\"\nReplace-ExactText $d \"63+28=91\" \"16+33=49\"\nReplace-ExactText $d \"44-40=4\" \"75-52=23\"\nReplace-ExactText $d \"32+29=61\" \"20+7=27\"\nReplace-ExactText $d \"27+25=52\" \"48+48=96\"\nReplace-ExactText $d \"55-38=17\" \"47-42=5\"\nReplace-ExactText $d \"6-3=3\" \"48-32=16\"\nReplace-ExactText $d \"38+26=64\" \"14+72=86\"\nReplace-ExactText $d \"44+35=79\" \"1+4=5\"\nReplace-ExactText $d \"32+54=86\" \"59-4=55\"\nReplace-ExactText $d \"7+52=59\" \"51-23=28\"\nReplace-ExactText $d \"34+64=98\" \"98-73=25\"\nReplace-ExactText $d \"78-43=35\" \"78-36=42\"\nReplace-ExactText $d \"82-36=46\" \"8-3=5\"\nReplace-ExactText $d \"55+4=59\" \"19+21=40\"\nReplace-ExactText $d \"62-27=35\" \"20+9=29\"\nReplace-ExactText $d \"74-3=71\" \"49+35=84\"\nReplace-ExactText $d \"93-80=13\" \"74-9=65\"\nReplace-ExactText $d \"37+49=86\" \"2+56=58\"\nReplace-ExactText $d \"47-13=34\" \"41-27=14\"\nReplace-ExactText $d \"92-1=91\" \"35+15=50\"\nReplace-ExactText $d \"93-16=77\" \"95-8=87\"\nReplace-ExactText $d \"0+79=79\" \"36+8=44\"\nReplace-ExactText $d \"27+61=88\" \"76+4=80\"\nReplace-ExactText $d \"24+5=29\" \"96-85=11\"\nReplace-ExactText $d \"93-65=28\" \"50+45=95\"\nReplace-ExactText $d \"58-52=6\" \"49+8=57\"\nReplace-ExactText $d \"12+54=66\" \"37-22=15\"\nReplace-ExactText $d \"67+15=82\" \"70-57=13\"\nReplace-ExactText $d \"79-67=12\" \"7+70=77\"\nReplace-ExactText $d \"88-53=35\" \"8+87=95\"\nReplace-ExactText $d \"13-10=3\" \"23+37=60\"\nReplace-ExactText $d \"30+29=59\" \"62+16=78\"\nReplace-ExactText $d \"14+21=35\" \"70+11=81\"\nReplace-ExactText $d \"2+3=5\" \"42-4=38\"\nReplace-ExactText $d \"47+7=54\" \"64-13=51\"\nReplace-ExactText $d \"52-6=46\" \"39-6=33\"\nReplace-ExactText $d \"54-53=1\" \"3-2=1\"\nReplace-ExactText $d \"59+15=74\" \"54-51=3\"\nReplace-ExactText $d \"37-13=24\" \"48+42=90\"\nReplace-ExactText $d \"51+11=62\" \"77+21=98\"\nReplace-ExactText $d \"90-81=9\" \"93-19=74\"\nReplace-ExactText $d \"62-33=29\" \"86-63=23\"\nReplace-ExactText $d \"51+39=90\" \"67-9=58\"\nReplace-ExactText $d \"55-47=8\" \"15+41=56\"\nReplace-ExactText $d \"16+49=65\" \"75+5=80\"\nReplace-ExactText $d \"63-48=15\" \"83+3=86\"\nReplace-ExactText $d \"70-29=41\" \"40+21=61\"\nReplace-ExactText $d \"48+26=74\" \"34-33=1\"\nReplace-ExactText $d \"52+39=91\" \"46-43=3\"\nReplace-ExactText $d \"31+40=71\" \"69-55=14\"\nReplace-ExactText $d \"34+46=80\" \"0+63=63\"\nReplace-ExactText $d \"11+1=12\" \"45+47=92\"\nReplace-ExactText $d \"80-12=68\" \"59-42=17\"\nReplace-ExactText $d \"39+20=59\" \"18+79=97\"\nReplace-ExactText $d \"8+28=36\" \"21+1=22\"\nReplace-ExactText $d \"63+1=64\" \"45-4=41\"\nReplace-ExactText $d \"23+18=41\" \"68-23=45\"\nReplace-ExactText $d \"3+88=91\" \"30+64=94\"\nReplace-ExactText $d \"28+52=80\" \"58-8=50\"\nReplace-ExactText $d \"25+62=87\" \"9+47=56\"\nReplace-ExactText $d \"94-75=19\" \"65-41=24\"\nReplace-ExactText $d \"55+33=88\" \"17+46=63\"\nReplace-ExactText $d \"30-26=4\" \"40+19=59\"\nReplace-ExactText $d \"43+18=61\" \"85+3=88\"\nReplace-ExactText $d \"3+34=37\" \"44+1=45\"\nReplace-ExactText $d \"73-46=27\" \"81-36=45\"\nReplace-ExactText $d \"6+42=48\" \"32+60=92\"\nReplace-ExactText $d \"44+23=67\" \"46+2=48\"\nReplace-ExactText $d \"17+59=76\" \"0+6=6\"\nReplace-ExactText $d \"27-19=8\" \"18+80=98\"\nReplace-ExactText $d \"89-71=18\" \"78-58=20\"\nReplace-ExactText $d \"68-28=40\" \"60+16=76\"\nReplace-ExactText $d \"85+6=91\" \"21+53=74\"\nReplace-ExactText $d \"69-39=30\" \"94-57=37\"\nReplace-ExactText $d \"34-29=5\" \"12+41=53\"\nReplace-ExactText $d \"3+92=95\" \"38+34=72\"\nReplace-ExactText $d \"89-70=19\" \"98-0=98\"\nReplace-ExactText $d \"82-4=78\" \"51+34=85\"\nReplace-ExactText $d \"75-45=30\" \"94-93=1\"\nReplace-ExactText $d \"66-49=17\" \"78-45=33\"\nReplace-ExactText $d \"34+21=55\" \"70-18=52\"\nReplace-ExactText $d \"30-29=1\" \"71-17=54\"\nReplace-ExactText $d \"13-2=11\" \"9+62=71\"\nReplace-ExactText $d \"38+50=88\" \"10-4=6\"\nReplace-ExactText $d \"34+26=60\" \"44+32=76\"\nReplace-ExactText $d \"54+40=94\" \"8+83=91\"\nReplace-ExactText $d \"79+20=99\" \"33+55=88\"\nReplace-ExactText $d \"85-32=53\" \"69-7=62\"\nReplace-ExactText $d \"96-86=10\" \"13+85=98\"\nReplace-ExactText $d \"91-70=21\" \"69-43=26\"\nReplace-ExactText $d \"10+38=48\" \"12+14=26\"\nReplace-ExactText $d \"71+19=90\" \"89-18=71\"\nReplace-ExactText $d \"46-42=4\" \"5+51=56\"\nReplace-ExactText $d \"81+8=89\" \"92-80=12\"\nReplace-ExactText $d \"73-5=68\" \"90+6=96\"\nReplace-ExactText $d \"88-30=58\" \"65+7=72\"\nReplace-ExactText $d \"34-31=3\" \"92-27=65\"\n"}
